$wb = $excel.ActiveWorkbook

# The "LoginTest" worksheet's expected-result cell (D2) changes from
# "Positive" to "Negative".
$wsLogin = $wb.Worksheets.Item("LoginTest")
$wsLogin.Range("D2").Value = "Negative"

# LoginTest becomes the active sheet/tab (previously "TestCases" was),
# with D3 as the selected cell.
$wsLogin.Activate()
$wsLogin.Range("D3").Select()
